# Generate Report for Handoff
# A new localization handoff run (455b784e-b379-49f0-8c0b-0d528e073873) has
# completed. Insert its results as the new "latest" row (row 2) on every
# sheet, pushing the previous latest run (5e8bf5fb-8e46-44e0-a234-e063797f8d42)
# down to row 3, and re-wire the hyperlinks to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Hyperlinks.Delete()
$ws.Rows("2:2").Insert()

$ws.Range("A2").Value2 = "455b784e-b379-49f0-8c0b-0d528e073873.md"
$ws.Range("B2").Value2 = "Ready for handoff"
$ws.Range("C2").Value2 = "Ready for handoff"
$ws.Range("D2").Value2 = "2016-26-11 20:26:05"

$ws.Range("A3").Value2 = "5e8bf5fb-8e46-44e0-a234-e063797f8d42.md"
$ws.Range("B3").Value2 = "Ready for handoff"
$ws.Range("C3").Value2 = "Ready for handoff"
$ws.Range("D3").Value2 = "2016-25-11 20:25:49"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/b6d0beef57d7aee95bb469f3eb15be071f1014e9/e2e/455b784e-b379-49f0-8c0b-0d528e073873.md", "", "", "455b784e-b379-49f0-8c0b-0d528e073873.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/b6d0beef57d7aee95bb469f3eb15be071f1014e9/e2e/5e8bf5fb-8e46-44e0-a234-e063797f8d42.md", "", "", "5e8bf5fb-8e46-44e0-a234-e063797f8d42.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Hyperlinks.Delete()
$ws.Rows("2:2").Insert()

$ws.Range("A2").Value2 = "455b784e-b379-49f0-8c0b-0d528e073873.md"
$ws.Range("B2").Value2 = ".md"
$ws.Range("C2").Value2 = "Ready for handoff"
$ws.Range("D2").Value2 = "455b784e-b379-49f0-8c0b-0d528e073873.89e979a510f7d4a5489f4a7ef207a128346d06f0.zh-cn.xlf"
$ws.Range("E2").Value2 = "2016-03-11 20:26:02"
$ws.Range("H2").Value2 = "0001-01-01 00:00:00"
$ws.Range("I2").Value2 = "Include"

$ws.Range("A3").Value2 = "5e8bf5fb-8e46-44e0-a234-e063797f8d42.md"
$ws.Range("B3").Value2 = ".md"
$ws.Range("C3").Value2 = "Ready for handoff"
$ws.Range("D3").Value2 = "5e8bf5fb-8e46-44e0-a234-e063797f8d42.03007089372f3fce5b254323b6067f9a52eeb39f.zh-cn.xlf"
$ws.Range("E3").Value2 = "2016-03-11 20:25:45"
$ws.Range("H3").Value2 = "0001-01-01 00:00:00"
$ws.Range("I3").Value2 = "Include"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/b6d0beef57d7aee95bb469f3eb15be071f1014e9/e2e/455b784e-b379-49f0-8c0b-0d528e073873.md", "", "", "455b784e-b379-49f0-8c0b-0d528e073873.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/b6d0beef57d7aee95bb469f3eb15be071f1014e9/e2e/455b784e-b379-49f0-8c0b-0d528e073873.md", "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c30e8d7f361e058613a52112fff8760775da1d98/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/455b784e-b379-49f0-8c0b-0d528e073873.89e979a510f7d4a5489f4a7ef207a128346d06f0.zh-cn.xlf", "", "", "455b784e-b379-49f0-8c0b-0d528e073873.89e979a510f7d4a5489f4a7ef207a128346d06f0.zh-cn.xlf") | Out-Null

$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/b6d0beef57d7aee95bb469f3eb15be071f1014e9/e2e/5e8bf5fb-8e46-44e0-a234-e063797f8d42.md", "", "", "5e8bf5fb-8e46-44e0-a234-e063797f8d42.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/b6d0beef57d7aee95bb469f3eb15be071f1014e9/e2e/5e8bf5fb-8e46-44e0-a234-e063797f8d42.md", "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c30e8d7f361e058613a52112fff8760775da1d98/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/5e8bf5fb-8e46-44e0-a234-e063797f8d42.03007089372f3fce5b254323b6067f9a52eeb39f.zh-cn.xlf", "", "", "5e8bf5fb-8e46-44e0-a234-e063797f8d42.03007089372f3fce5b254323b6067f9a52eeb39f.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Hyperlinks.Delete()
$ws.Rows("2:2").Insert()

$ws.Range("A2").Value2 = "455b784e-b379-49f0-8c0b-0d528e073873.md"
$ws.Range("B2").Value2 = ".md"
$ws.Range("C2").Value2 = "Ready for handoff"
$ws.Range("D2").Value2 = "455b784e-b379-49f0-8c0b-0d528e073873.89e979a510f7d4a5489f4a7ef207a128346d06f0.de-de.xlf"
$ws.Range("E2").Value2 = "2016-03-11 20:26:05"
$ws.Range("H2").Value2 = "0001-01-01 00:00:00"
$ws.Range("I2").Value2 = "Include"

$ws.Range("A3").Value2 = "5e8bf5fb-8e46-44e0-a234-e063797f8d42.md"
$ws.Range("B3").Value2 = ".md"
$ws.Range("C3").Value2 = "Ready for handoff"
$ws.Range("D3").Value2 = "5e8bf5fb-8e46-44e0-a234-e063797f8d42.03007089372f3fce5b254323b6067f9a52eeb39f.de-de.xlf"
$ws.Range("E3").Value2 = "2016-03-11 20:25:49"
$ws.Range("H3").Value2 = "0001-01-01 00:00:00"
$ws.Range("I3").Value2 = "Include"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/b6d0beef57d7aee95bb469f3eb15be071f1014e9/e2e/455b784e-b379-49f0-8c0b-0d528e073873.md", "", "", "455b784e-b379-49f0-8c0b-0d528e073873.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/b6d0beef57d7aee95bb469f3eb15be071f1014e9/e2e/455b784e-b379-49f0-8c0b-0d528e073873.md", "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8a0567468a383bb0cc360315d47588f1c222d8cc/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/455b784e-b379-49f0-8c0b-0d528e073873.89e979a510f7d4a5489f4a7ef207a128346d06f0.de-de.xlf", "", "", "455b784e-b379-49f0-8c0b-0d528e073873.89e979a510f7d4a5489f4a7ef207a128346d06f0.de-de.xlf") | Out-Null

$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/b6d0beef57d7aee95bb469f3eb15be071f1014e9/e2e/5e8bf5fb-8e46-44e0-a234-e063797f8d42.md", "", "", "5e8bf5fb-8e46-44e0-a234-e063797f8d42.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/b6d0beef57d7aee95bb469f3eb15be071f1014e9/e2e/5e8bf5fb-8e46-44e0-a234-e063797f8d42.md", "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8a0567468a383bb0cc360315d47588f1c222d8cc/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/5e8bf5fb-8e46-44e0-a234-e063797f8d42.03007089372f3fce5b254323b6067f9a52eeb39f.de-de.xlf", "", "", "5e8bf5fb-8e46-44e0-a234-e063797f8d42.03007089372f3fce5b254323b6067f9a52eeb39f.de-de.xlf") | Out-Null

$wb.Worksheets.Item("Overview").Activate()
